$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference style (style index 0 / General, no borders) used to restore
# cell formatting after forcing text entry for numeric-looking values.
$normalStyle = $ws.Range("C2").Style

# Row 2
$ws.Range("D2").Value = "67.940.81"
$ws.Range("E2").Value = "  +0.02%  "

# Row 3
$ws.Range("D3").Value = "3.251.29"
$ws.Range("E3").Value = "  -0.79%  "

# Row 4
$ws.Range("E4").Value = "  +0.09%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "581.67"
$ws.Range("D5").Style = $normalStyle
$ws.Range("E5").Value = "  +0.08%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "184.15"
$ws.Range("D6").Style = $normalStyle
$ws.Range("E6").Value = "  +0.57%  "

# Row 7
$ws.Range("E7").Value = "  +0.06%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.597"
$ws.Range("D8").Style = $normalStyle
$ws.Range("E8").Value = "  -1.13%  "

# Row 9
$ws.Range("E9").Value = "  -1.99%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.64"
$ws.Range("D10").Style = $normalStyle
$ws.Range("E10").Value = "  -1.25%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.418"
$ws.Range("D11").Style = $normalStyle
$ws.Range("E11").Value = "  +0.08%  "

# Row 12
$ws.Range("D12").Value = "3.820.53"
$ws.Range("E12").Value = "  -0.72%  "

# Row 13
$ws.Range("E13").Value = "  -0.17%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.08"
$ws.Range("D14").Style = $normalStyle
$ws.Range("E14").Value = "  -2.23%  "

# Row 15
$ws.Range("D15").Value = "68.147.24"
$ws.Range("E15").Value = "  +0.41%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000170"
$ws.Range("D16").Style = $normalStyle
$ws.Range("E16").Value = "  +0.41%  "

# Row 17
$ws.Range("D17").Value = "3.235.12"
$ws.Range("E17").Value = "  -1.15%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.83"
$ws.Range("D18").Style = $normalStyle
$ws.Range("E18").Value = "  -0.66%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.57"
$ws.Range("D19").Style = $normalStyle
$ws.Range("E19").Value = "  +0.23%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "392.31"
$ws.Range("D20").Style = $normalStyle
$ws.Range("E20").Value = "  +3.87%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.68"
$ws.Range("D21").Style = $normalStyle
$ws.Range("E21").Value = "  -0.01%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "71.43"
$ws.Range("D22").Style = $normalStyle
$ws.Range("E22").Value = "  +0.11%  "

# Row 23
$ws.Range("E23").Value = "  -0.13%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.518"
$ws.Range("D24").Style = $normalStyle
$ws.Range("E24").Value = "  +0.56%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000119"
$ws.Range("D25").Style = $normalStyle
$ws.Range("E25").Value = "  -1.05%  "

# Row 26
$ws.Range("E26").Value = "  +4.24%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.70"
$ws.Range("D27").Style = $normalStyle
$ws.Range("E27").Value = "  -0.20%  "

# Row 28
$ws.Range("E28").Value = "  +0.04%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.97"
$ws.Range("D29").Style = $normalStyle
$ws.Range("E29").Value = "  -0.53%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.68"
$ws.Range("D30").Style = $normalStyle
$ws.Range("E30").Value = "  -0.68%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "22.86"
$ws.Range("D31").Style = $normalStyle
$ws.Range("E31").Value = "  -0.33%  "

# Row 32
$ws.Range("E32").Value = "  +2.41%  "

# Row 33
$ws.Range("B33").Value = "USDe"
$ws.Range("C33").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.998"
$ws.Range("D33").Style = $normalStyle
$ws.Range("E33").Value = "  -0.07%  "

# Row 34
$ws.Range("B34").Value = "Fetch.AI"
$ws.Range("C34").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.28"
$ws.Range("D34").Style = $normalStyle
$ws.Range("E34").Value = "  -0.25%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "162.70"
$ws.Range("D35").Style = $normalStyle
$ws.Range("E35").Value = "  +0.60%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.49"
$ws.Range("D36").Style = $normalStyle
$ws.Range("E36").Value = "  -3.01%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.95"
$ws.Range("D37").Style = $normalStyle
$ws.Range("E37").Value = "  +4.97%  "

# Row 38
$ws.Range("B38").Value = "EnergySwap"
$ws.Range("C38").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "26.77"
$ws.Range("D38").Style = $normalStyle
$ws.Range("E38").Value = "  -0.97%  "

# Row 39
$ws.Range("B39").Value = "Mantle"
$ws.Range("C39").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.819"
$ws.Range("D39").Style = $normalStyle
$ws.Range("E39").Value = "  -4.09%  "

# Row 40
$ws.Range("E40").Value = "  -1.28%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.51"
$ws.Range("D41").Style = $normalStyle
$ws.Range("E41").Value = "  -4.05%  "

# Row 42
$ws.Range("E42").Value = "  -6.87%  "

# Row 43
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "2.653.51"
$ws.Range("E43").Value = "  +0.01%  "

# Row 44
$ws.Range("B44").Value = "Hedera"
$ws.Range("C44").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0688"
$ws.Range("D44").Style = $normalStyle
$ws.Range("E44").Value = "  +0.80%  "

# Row 45
$ws.Range("B45").Value = "InjectiveProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "25.37"
$ws.Range("D45").Style = $normalStyle
$ws.Range("E45").Value = "  -1.72%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "41.07"
$ws.Range("D46").Style = $normalStyle
$ws.Range("E46").Value = "  +0.19%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "338.26"
$ws.Range("D47").Style = $normalStyle
$ws.Range("E47").Value = "  -3.77%  "

# Row 48
$ws.Range("E48").Value = "  -1.26%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.35"
$ws.Range("D49").Style = $normalStyle
$ws.Range("E49").Value = "  +2.88%  "

# Row 50
$ws.Range("B50").Value = "Stellar"
$ws.Range("C50").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.102"
$ws.Range("D50").Style = $normalStyle
$ws.Range("E50").Value = "  -1.50%  "

# Row 51
$ws.Range("B51").Value = "Arweave"
$ws.Range("C51").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "31.36"
$ws.Range("D51").Style = $normalStyle
$ws.Range("E51").Value = "  +0.83%  "
